# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Espárragos".
# The new record is inserted at row 23, pushing the existing rows 23-57 down to 24-58
# (dimension grows from A1:R57 to A1:R58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 23..57 down by one to make room for the new record at row 23.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new observation.
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44874
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 300000000
$ws.Range("G23").Value = "Espárragos"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1700
$ws.Range("M23").Value = 1600
$ws.Range("N23").Value = "$/kilo"
$ws.Range("O23").Value = "Provincia de Linares"
$ws.Range("P23").Value = 1600
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
